$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.226.41'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +5.55%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.998.87'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.00%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.12'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.09%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '162.96'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +12.53%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('E8').Value = '  +3.06%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.993.70'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.88%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.52'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.39%  '

# Row 11
$ws.Range('E11').Value = '  +3.28%  '

# Row 12
$ws.Range('E12').Value = '  +4.67%  '

# Row 13
$ws.Range('E13').Value = '  +5.61%  '

# Row 14
$ws.Range('E14').Value = '  +4.61%  '

# Row 15
$ws.Range('E15').Value = '  -0.88%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.191.71'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +5.36%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.497.27'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.03%  '

# Row 18
$ws.Range('E18').Value = '  +3.88%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.001.38'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.93%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '454.50'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.27%  '

# Row 21
$ws.Range('E21').Value = '  +5.11%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.686'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.62%  '

# Row 23
$ws.Range('E23').Value = '  +5.88%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.29'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.24%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.31'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +13.45%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.24'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.01%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.38'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.28%  '

# Row 28
$ws.Range('E28').Value = '  +0.02%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.14'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +12.76%  '

# Row 30
$ws.Range('E30').Value = '  +18.01%  '

# Row 31
$ws.Range('E31').Value = '  +4.85%  '

# Row 32
$ws.Range('E32').Value = '  -5.27%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.29'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.19%  '

# Row 34
$ws.Range('E34').Value = '  +3.12%  '

# Row 35
$ws.Range('E35').Value = '  -0.18%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.992'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.59%  '

# Row 37
$ws.Range('E37').Value = '  +7.33%  '

# Row 38
$ws.Range('E38').Value = '  +7.78%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '49.68'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.15%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.97'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.20%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.310'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +15.42%  '

# Row 42
$ws.Range('E42').Value = '  +6.77%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '43.84'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +7.13%  '

# Row 44
$ws.Range('E44').Value = '  +3.60%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '399.31'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +12.11%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0358'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.21%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.787.32'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.61%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '133.51'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.05%  '

# Row 49
$ws.Range('E49').Value = '  +0.01%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.83'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +10.62%  '

# Row 51
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.16'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +10.10%  '
